$wb = $excel.ActiveWorkbook

# --- Sheet 'data' (sheet1): rows 37-66 ---
$wsData = $wb.Worksheets.Item("data")

# Row 37
$wsData.Cells.Item(37, 4).Value = 0.38
$wsData.Cells.Item(37, 5).Value = 0.25
$wsData.Cells.Item(37, 6).Value = 0.18
$wsData.Cells.Item(37, 7).Value = 0.19
$wsData.Cells.Item(37, 8).Value = 0.19
$wsData.Cells.Item(37, 9).Value = 0.25

# Row 38
$wsData.Cells.Item(38, 4).Value = 0.22
$wsData.Cells.Item(38, 5).Value = 0.2
$wsData.Cells.Item(38, 6).Value = 0.18
$wsData.Cells.Item(38, 8).Value = 0.15
$wsData.Cells.Item(38, 9).Value = 0.15

# Row 39
$wsData.Cells.Item(39, 4).Value = 0.1
$wsData.Cells.Item(39, 5).Value = 0.14
$wsData.Cells.Item(39, 6).Value = 0.12
$wsData.Cells.Item(39, 7).Value = 0.21
$wsData.Cells.Item(39, 8).Value = 0.21

# Row 40
$wsData.Cells.Item(40, 4).Value = 0.1
$wsData.Cells.Item(40, 6).Value = 0.23
$wsData.Cells.Item(40, 7).Value = 0.16
$wsData.Cells.Item(40, 8).Value = 0.18
$wsData.Cells.Item(40, 9).Value = 0.16

# Row 41
$wsData.Cells.Item(41, 4).Value = 0.2
$wsData.Cells.Item(41, 5).Value = 0.25
$wsData.Cells.Item(41, 6).Value = 0.29
$wsData.Cells.Item(41, 7).Value = 0.3
$wsData.Cells.Item(41, 8).Value = 0.27
$wsData.Cells.Item(41, 9).Value = 0.21

# Row 42
$wsData.Cells.Item(42, 4).Value = 0.3
$wsData.Cells.Item(42, 5).Value = 0.24
$wsData.Cells.Item(42, 6).Value = 0.21
$wsData.Cells.Item(42, 7).Value = 0.25
$wsData.Cells.Item(42, 8).Value = 0.24
$wsData.Cells.Item(42, 9).Value = 0.31

# Row 43
$wsData.Cells.Item(43, 4).Value = 0.22
$wsData.Cells.Item(43, 5).Value = 0.2
$wsData.Cells.Item(43, 6).Value = 0.25
$wsData.Cells.Item(43, 7).Value = 0.21
$wsData.Cells.Item(43, 8).Value = 0.23

# Row 44
$wsData.Cells.Item(44, 4).Value = 0.09
$wsData.Cells.Item(44, 5).Value = 0.09
$wsData.Cells.Item(44, 7).Value = 0.1
$wsData.Cells.Item(44, 8).Value = 0.16
$wsData.Cells.Item(44, 9).Value = 0.14

# Row 45
$wsData.Cells.Item(45, 4).Value = 0.14
$wsData.Cells.Item(45, 5).Value = 0.2
$wsData.Cells.Item(45, 6).Value = 0.27
$wsData.Cells.Item(45, 7).Value = 0.26
$wsData.Cells.Item(45, 8).Value = 0.18
$wsData.Cells.Item(45, 9).Value = 0.17

# Row 46
$wsData.Cells.Item(46, 4).Value = 0.25
$wsData.Cells.Item(46, 5).Value = 0.27
$wsData.Cells.Item(46, 6).Value = 0.16
$wsData.Cells.Item(46, 7).Value = 0.18
$wsData.Cells.Item(46, 8).Value = 0.19
$wsData.Cells.Item(46, 9).Value = 0.15

# Row 47
$wsData.Cells.Item(47, 4).Value = 0.28
$wsData.Cells.Item(47, 5).Value = 0.28
$wsData.Cells.Item(47, 6).Value = 0.35
$wsData.Cells.Item(47, 7).Value = 0.35
$wsData.Cells.Item(47, 8).Value = 0.36
$wsData.Cells.Item(47, 9).Value = 0.4

# Row 48
$wsData.Cells.Item(48, 4).Value = 0.25
$wsData.Cells.Item(48, 5).Value = 0.28
$wsData.Cells.Item(48, 6).Value = 0.25
$wsData.Cells.Item(48, 7).Value = 0.33
$wsData.Cells.Item(48, 8).Value = 0.33

# Row 49
$wsData.Cells.Item(49, 4).Value = 0.11
$wsData.Cells.Item(49, 5).Value = 0.11
$wsData.Cells.Item(49, 6).Value = 0.08
$wsData.Cells.Item(49, 7).Value = 0.1
$wsData.Cells.Item(49, 9).Value = 0.08

# Row 50
$wsData.Cells.Item(50, 4).Value = 0.17
$wsData.Cells.Item(50, 5).Value = 0.16
$wsData.Cells.Item(50, 6).Value = 0.2
$wsData.Cells.Item(50, 7).Value = 0.13
$wsData.Cells.Item(50, 8).Value = 0.16
$wsData.Cells.Item(50, 9).Value = 0.14

# Row 51
$wsData.Cells.Item(51, 4).Value = 0.19
$wsData.Cells.Item(51, 6).Value = 0.12

# Row 52
$wsData.Cells.Item(52, 4).Value = 0.26
$wsData.Cells.Item(52, 5).Value = 0.22
$wsData.Cells.Item(52, 6).Value = 0.25
$wsData.Cells.Item(52, 7).Value = 0.29
$wsData.Cells.Item(52, 8).Value = 0.27
$wsData.Cells.Item(52, 9).Value = 0.37

# Row 53
$wsData.Cells.Item(53, 4).Value = 0.25
$wsData.Cells.Item(53, 5).Value = 0.23
$wsData.Cells.Item(53, 6).Value = 0.21
$wsData.Cells.Item(53, 7).Value = 0.27
$wsData.Cells.Item(53, 8).Value = 0.23
$wsData.Cells.Item(53, 9).Value = 0.23

# Row 54
$wsData.Cells.Item(54, 4).Value = 0.14
$wsData.Cells.Item(54, 6).Value = 0.25
$wsData.Cells.Item(54, 7).Value = 0.18
$wsData.Cells.Item(54, 9).Value = 0.18

# Row 55
$wsData.Cells.Item(55, 4).Value = 0.13
$wsData.Cells.Item(55, 5).Value = 0.14
$wsData.Cells.Item(55, 6).Value = 0.16
$wsData.Cells.Item(55, 7).Value = 0.12
$wsData.Cells.Item(55, 8).Value = 0.15
$wsData.Cells.Item(55, 9).Value = 0.1

# Row 56
$wsData.Cells.Item(56, 4).Value = 0.22
$wsData.Cells.Item(56, 5).Value = 0.25
$wsData.Cells.Item(56, 6).Value = 0.13
$wsData.Cells.Item(56, 7).Value = 0.14
$wsData.Cells.Item(56, 8).Value = 0.14
$wsData.Cells.Item(56, 9).Value = 0.12

# Row 57
$wsData.Cells.Item(57, 4).Value = 0.39
$wsData.Cells.Item(57, 5).Value = 0.34
$wsData.Cells.Item(57, 6).Value = 0.4
$wsData.Cells.Item(57, 7).Value = 0.5
$wsData.Cells.Item(57, 8).Value = 0.48
$wsData.Cells.Item(57, 9).Value = 0.5600000000000001

# Row 58
$wsData.Cells.Item(58, 4).Value = 0.26
$wsData.Cells.Item(58, 5).Value = 0.29
$wsData.Cells.Item(58, 6).Value = 0.32
$wsData.Cells.Item(58, 7).Value = 0.28
$wsData.Cells.Item(58, 8).Value = 0.26
$wsData.Cells.Item(58, 9).Value = 0.24

# Row 59
$wsData.Cells.Item(59, 4).Value = 0.12
$wsData.Cells.Item(59, 6).Value = 0.12
$wsData.Cells.Item(59, 7).Value = 0.09
$wsData.Cells.Item(59, 8).Value = 0.11
$wsData.Cells.Item(59, 9).Value = 0.09

# Row 60
$wsData.Cells.Item(60, 4).Value = 0.06
$wsData.Cells.Item(60, 5).Value = 0.17
$wsData.Cells.Item(60, 6).Value = 0.11
$wsData.Cells.Item(60, 7).Value = 0.09
$wsData.Cells.Item(60, 8).Value = 0.1
$wsData.Cells.Item(60, 9).Value = 0.08

# Row 61
$wsData.Cells.Item(61, 4).Value = 0.17
$wsData.Cells.Item(61, 5).Value = 0.12
$wsData.Cells.Item(61, 6).Value = 0.05
$wsData.Cells.Item(61, 7).Value = 0.04
$wsData.Cells.Item(61, 8).Value = 0.05
$wsData.Cells.Item(61, 9).Value = 0.03

# Row 62
$wsData.Cells.Item(62, 4).Value = 0.48
$wsData.Cells.Item(62, 5).Value = 0.45
$wsData.Cells.Item(62, 6).Value = 0.55
$wsData.Cells.Item(62, 7).Value = 0.58
$wsData.Cells.Item(62, 8).Value = 0.63
$wsData.Cells.Item(62, 9).Value = 0.6899999999999999

# Row 63
$wsData.Cells.Item(63, 4).Value = 0.23
$wsData.Cells.Item(63, 5).Value = 0.27
$wsData.Cells.Item(63, 6).Value = 0.25
$wsData.Cells.Item(63, 7).Value = 0.2
$wsData.Cells.Item(63, 8).Value = 0.18
$wsData.Cells.Item(63, 9).Value = 0.16

# Row 64
$wsData.Cells.Item(64, 4).Value = 0.07000000000000001
$wsData.Cells.Item(64, 5).Value = 0.07000000000000001
$wsData.Cells.Item(64, 6).Value = 0.05
$wsData.Cells.Item(64, 8).Value = 0.06
$wsData.Cells.Item(64, 9).Value = 0.05

# Row 65
$wsData.Cells.Item(65, 5).Value = 0.13
$wsData.Cells.Item(65, 6).Value = 0.11
$wsData.Cells.Item(65, 7).Value = 0.09
$wsData.Cells.Item(65, 8).Value = 0.09
$wsData.Cells.Item(65, 9).Value = 0.07000000000000001

# Row 66
$wsData.Cells.Item(66, 4).Value = 0.09
$wsData.Cells.Item(66, 5).Value = 0.08
$wsData.Cells.Item(66, 6).Value = 0.04
$wsData.Cells.Item(66, 7).Value = 0.05
$wsData.Cells.Item(66, 8).Value = 0.04
$wsData.Cells.Item(66, 9).Value = 0.03


# --- Sheet 'pocetR' (sheet2): rows 9-14 ---
$wsPocetR = $wb.Worksheets.Item("pocetR")

# Row 9
$wsPocetR.Cells.Item(9, 3).Value = 504
$wsPocetR.Cells.Item(9, 4).Value = 494
$wsPocetR.Cells.Item(9, 5).Value = 520
$wsPocetR.Cells.Item(9, 6).Value = 500
$wsPocetR.Cells.Item(9, 7).Value = 504
$wsPocetR.Cells.Item(9, 8).Value = 511

# Row 10
$wsPocetR.Cells.Item(10, 3).Value = 510
$wsPocetR.Cells.Item(10, 4).Value = 511
$wsPocetR.Cells.Item(10, 5).Value = 511
$wsPocetR.Cells.Item(10, 6).Value = 497
$wsPocetR.Cells.Item(10, 7).Value = 500
$wsPocetR.Cells.Item(10, 8).Value = 492

# Row 11
$wsPocetR.Cells.Item(11, 3).Value = 294
$wsPocetR.Cells.Item(11, 4).Value = 288
$wsPocetR.Cells.Item(11, 5).Value = 295
$wsPocetR.Cells.Item(11, 6).Value = 285
$wsPocetR.Cells.Item(11, 7).Value = 281
$wsPocetR.Cells.Item(11, 8).Value = 284

# Row 12
$wsPocetR.Cells.Item(12, 3).Value = 486
$wsPocetR.Cells.Item(12, 4).Value = 492
$wsPocetR.Cells.Item(12, 5).Value = 489
$wsPocetR.Cells.Item(12, 6).Value = 484
$wsPocetR.Cells.Item(12, 7).Value = 477
$wsPocetR.Cells.Item(12, 8).Value = 481

# Row 13
$wsPocetR.Cells.Item(13, 3).Value = 234
$wsPocetR.Cells.Item(13, 4).Value = 233
$wsPocetR.Cells.Item(13, 5).Value = 232
$wsPocetR.Cells.Item(13, 6).Value = 229
$wsPocetR.Cells.Item(13, 7).Value = 223
$wsPocetR.Cells.Item(13, 8).Value = 227

# Row 14
$wsPocetR.Cells.Item(14, 3).Value = 139
$wsPocetR.Cells.Item(14, 4).Value = 137
$wsPocetR.Cells.Item(14, 5).Value = 139
$wsPocetR.Cells.Item(14, 6).Value = 136
$wsPocetR.Cells.Item(14, 7).Value = 135
$wsPocetR.Cells.Item(14, 8).Value = 135

